# "Generate Report for Archive"
# - Update the localization "Status" text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview summary columns for
#   zh-cn/de-de, and the Status column on each language detail sheet).
# - Narrow the now-shorter Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 12.5   # characters; renders to the narrowest reachable column width

# --- Overview sheet: columns E (zh-cn status) and F (de-de status) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: column C (Status) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: column C (Status) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
